$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74
$values = @("2024-09-25T18:06:40Z", "temperature", "25", "N/A", "N/A", "N/A")

for ($col = 1; $col -le 6; $col++) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $values[$col - 1]
    $c.ClearFormats()
}
